$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 to the new
# header cells I1 and J1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for rows 2-14 (I0 and IF columns hold the same values)
$values = @(
    @(2, 8),
    @(3, 6),
    @(4, 9),
    @(5, 6),
    @(6, 9),
    @(7, 3),
    @(8, 9),
    @(9, 9),
    @(10, 4),
    @(11, 5),
    @(12, 5),
    @(13, 3),
    @(14, 3)
)

foreach ($row in $values) {
    $r = $row[0]
    $v = $row[1]
    $ws.Cells.Item($r, 9).Value = $v
    $ws.Cells.Item($r, 10).Value = $v
}
